$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D as text first, so numeric-looking price strings
# (e.g. "191.25") are not silently converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Price (column D) updates ---
$ws.Range("D2").Value = "69.671.70"
$ws.Range("D3").Value = "3.389.69"
$ws.Range("D5").Value = "191.25"
$ws.Range("D6").Value = "594.08"
$ws.Range("D8").Value = "0.608"
$ws.Range("D9").Value = "0.134"
$ws.Range("D10").Value = "6.77"
$ws.Range("D11").Value = "0.419"
$ws.Range("D12").Value = "3.981.05"
$ws.Range("D14").Value = "28.74"
$ws.Range("D15").Value = "69.633.41"
$ws.Range("D17").Value = "3.387.95"
$ws.Range("D18").Value = "451.30"
$ws.Range("D19").Value = "5.84"
$ws.Range("D20").Value = "13.83"
$ws.Range("D21").Value = "7.83"
$ws.Range("D22").Value = "76.37"
$ws.Range("D24").Value = "0.523"
$ws.Range("D27").Value = "9.50"
$ws.Range("D28").Value = "1.00"
$ws.Range("D29").Value = "2.02"
$ws.Range("D30").Value = "23.49"
$ws.Range("D31").Value = "5.65"
$ws.Range("D32").Value = "1.28"
$ws.Range("D33").Value = "7.00"
$ws.Range("D35").Value = "1.57"
$ws.Range("D36").Value = "165.54"
$ws.Range("D38").Value = "28.36"
$ws.Range("D39").Value = "0.816"
$ws.Range("D40").Value = "4.61"
$ws.Range("D42").Value = "2.759.65"
$ws.Range("D43").Value = "2.52"
$ws.Range("D44").Value = "25.60"
$ws.Range("D45").Value = "0.0689"
$ws.Range("D46").Value = "41.12"
$ws.Range("D47").Value = "341.05"
$ws.Range("D48").Value = "0.0285"
$ws.Range("D49").Value = "33.12"
$ws.Range("D51").Value = "6.34"

# Restore column D style back to the default/normal style so the
# underlying cell XML does not carry a leftover numFmt style index.
$ws.Range("D2:D51").Style = "Normal"

# --- Volume(1h) (column E) updates ---
$ws.Range("E3").Value = "  +4.30%  "
$ws.Range("E5").Value = "  +3.93%  "
$ws.Range("E6").Value = "  +2.42%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("E9").Value = "  +2.55%  "
$ws.Range("E10").Value = "  +2.95%  "
$ws.Range("E11").Value = "  +1.95%  "
$ws.Range("E12").Value = "  +4.78%  "
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("E14").Value = "  +3.86%  "
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("E17").Value = "  +5.44%  "
$ws.Range("E18").Value = "  +14.37%  "
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("E20").Value = "  +2.23%  "
$ws.Range("E21").Value = "  +3.31%  "
$ws.Range("E22").Value = "  +6.86%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("E25").Value = "  +4.28%  "
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E29").Value = "  +2.95%  "
$ws.Range("E30").Value = "  +3.76%  "
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("E32").Value = "  +2.44%  "
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  +6.55%  "
$ws.Range("E36").Value = "  +2.47%  "
$ws.Range("E37").Value = "  +2.66%  "
$ws.Range("E38").Value = "  +6.53%  "
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("E40").Value = "  +1.58%  "
$ws.Range("E41").Value = "  +2.17%  "
$ws.Range("E42").Value = "  +5.44%  "
$ws.Range("E43").Value = "  +2.02%  "
$ws.Range("E44").Value = "  +3.53%  "
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("E47").Value = "  +1.97%  "
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("E49").Value = "  +7.70%  "
$ws.Range("E50").Value = "  +5.50%  "
$ws.Range("E51").Value = "  +0.35%  "
